$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 9
$ws.Range("B4").Value = 0.7
